# KENTUCKY_2019.xlsx network-data cleanup
#
# 1) Rename the header row (A1:D1) from the Spanish labels to the short
#    machine-friendly column names used by the cleaning scripts.
# 2) Proper-case (title-case) the Mexican state/municipality names in
#    columns A and B for every data row, so connector words like
#    "de"/"del"/"de la"/"de los"/"el"/"y" are capitalised consistently
#    with the rest of the words (matches Excel's PROPER() behaviour for
#    this data, which only ever splits on spaces here).
# 3) Drop the trailing sample-size / source / credits rows that used to
#    sit below the data (rows 861-866), shrinking the used range back
#    down to A1:D860.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header rename -------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2) Title-case columns A and B for the data rows -------------------
for ($r = 2; $r -le 860; $r++) {

    # NOTE: this host's -eq/-ne/-ceq/-cne all compare case-insensitively,
    # so we can't guard the write with "if changed". Just write the
    # recomputed (title-cased) value back unconditionally - idempotent
    # for cells that were already in the target case.
    #
    # NOTE 2: plain "+" between two strings that both happen to look like
    # numbers (e.g. "2" + "2") is evaluated as NUMERIC addition by this
    # host ("2"+"2" -> 4, not "22"), which silently corrupts names like
    # "San Pedro Mixtepec - Distr. 22 -". Use [string]::Concat(...) for
    # per-word capitalisation and the -join operator to rebuild the
    # sentence so everything stays a string.

    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value()
    if ($aVal -ne $null -and $aVal -is [string] -and $aVal.Length -gt 0) {
        $words = $aVal.Split(" ")
        for ($i = 0; $i -lt $words.Length; $i++) {
            $w = $words[$i]
            if ($w.Length -gt 0) {
                $words[$i] = [string]::Concat($w.Substring(0,1).ToUpper(), $w.Substring(1))
            }
        }
        $aCell.Value = ($words -join " ")
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    if ($bVal -ne $null -and $bVal -is [string] -and $bVal.Length -gt 0) {
        $words = $bVal.Split(" ")
        for ($i = 0; $i -lt $words.Length; $i++) {
            $w = $words[$i]
            if ($w.Length -gt 0) {
                $words[$i] = [string]::Concat($w.Substring(0,1).ToUpper(), $w.Substring(1))
            }
        }
        $bCell.Value = ($words -join " ")
    }
}

# --- 3) Remove the trailing metadata rows below the data ---------------
$ws.Rows("861:866").Delete()

# --- 4) Re-derived percentage for the San Luis Potosi subtotal ---------
# The upstream cleaning script recomputes this percentage slightly
# differently than the original export, landing one ULP away
# (0.009556907037358819 -> ...821) even though both equal 33/3453.
$ws.Range("D629").Value = 0.009556907037358821
